# Swap the presentation's applied theme colour palette from the
# "Integral" (Red Violet) scheme to the default Office "Office Theme"
# colour scheme.
#
# ppt/theme/theme1.xml (the theme used by the slide master / all
# slides) keeps the same 12 colour slots (dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink) but each slot's RGB value is changed
# from the Integral theme's value to the standard Office theme value,
# matching the target diff.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.ColorScheme

# Index -> (scheme slot, target "Office Theme" RGB hex)
# Order matches PowerPoint's ColorScheme.Colors() enumeration:
#   1=dk1 2=lt1 3=dk2 4=lt2 5=accent1 6=accent2 7=accent3 8=accent4
#   9=accent5 10=accent6 11=hlink 12=folHlink
$officeThemeRGB = @(
    0,         # dk1      000000
    16777215,  # lt1      FFFFFF
    6968388,   # dk2      44546A
    15132391,  # lt2      E7E6E6
    13998939,  # accent1  5B9BD5
    3243501,   # accent2  ED7D31
    10855845,  # accent3  A5A5A5
    49407,     # accent4  FFC000
    12874308,  # accent5  4472C4
    4697456,   # accent6  70AD47
    12673797,  # hlink    0563C1
    7491477    # folHlink 954F72
)

for ($i = 1; $i -le $officeThemeRGB.Count; $i++) {
    $colorScheme.Colors($i).RGB = $officeThemeRGB[$i - 1]
}
